$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 1467
$ws.Cells.Item(3, 6).Value = 1440
$ws.Cells.Item(4, 6).Value = 403
$ws.Cells.Item(5, 6).Value = 225
$ws.Cells.Item(6, 6).Value = 699
$ws.Cells.Item(8, 6).Value = 632
$ws.Cells.Item(11, 6).Value = 1381
$ws.Cells.Item(12, 6).Value = 32955
$ws.Cells.Item(13, 6).Value = 7064
$ws.Cells.Item(14, 6).Value = 110
$ws.Cells.Item(15, 6).Value = 365
$ws.Cells.Item(16, 6).Value = 572
$ws.Cells.Item(17, 6).Value = 440
$ws.Cells.Item(19, 6).Value = 105
$ws.Cells.Item(20, 6).Value = 17
$ws.Cells.Item(22, 6).Value = 448
$ws.Cells.Item(23, 6).Value = 105
$ws.Cells.Item(24, 6).Value = 800
$ws.Cells.Item(25, 6).Value = 8
$ws.Cells.Item(26, 6).Value = 318
$ws.Cells.Item(27, 6).Value = 388
$ws.Cells.Item(28, 6).Value = 441
$ws.Cells.Item(30, 6).Value = 209
$ws.Cells.Item(32, 6).Value = 739
$ws.Cells.Item(35, 6).Value = 739

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 1187
$ws.Cells.Item(5, 6).Value = 161
$ws.Cells.Item(9, 6).Value = 234
$ws.Cells.Item(19, 6).Value = 4296

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 1453
$ws.Cells.Item(3, 6).Value = 354

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 1453
$ws.Cells.Item(3, 6).Value = 354
$ws.Cells.Item(4, 6).Value = 1187
$ws.Cells.Item(5, 6).Value = 1467
$ws.Cells.Item(7, 6).Value = 1440
$ws.Cells.Item(8, 6).Value = 225
$ws.Cells.Item(9, 6).Value = 699
$ws.Cells.Item(11, 6).Value = 632
$ws.Cells.Item(13, 6).Value = 1381
$ws.Cells.Item(14, 6).Value = 161
$ws.Cells.Item(17, 6).Value = 234
$ws.Cells.Item(18, 6).Value = 234
$ws.Cells.Item(21, 6).Value = 7064
$ws.Cells.Item(22, 6).Value = 110
$ws.Cells.Item(23, 6).Value = 365
$ws.Cells.Item(25, 6).Value = 572
$ws.Cells.Item(26, 6).Value = 440
$ws.Cells.Item(28, 6).Value = 105
$ws.Cells.Item(31, 6).Value = 448
$ws.Cells.Item(32, 6).Value = 105
$ws.Cells.Item(33, 6).Value = 800
$ws.Cells.Item(34, 6).Value = 318
$ws.Cells.Item(35, 6).Value = 388
$ws.Cells.Item(36, 6).Value = 441
$ws.Cells.Item(38, 6).Value = 209
$ws.Cells.Item(40, 6).Value = 739
